$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6756153702735901
$ws.Range("B1").Value = 2.104029178619385
$ws.Range("C1").Value = 5.113800525665283
$ws.Range("D1").Value = 2.934109926223755
$ws.Range("E1").Value = 0.6591125130653381
